$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel (this engine) implements Range.Hyperlinks.Delete() as "wipe every
# hyperlink on the sheet", so do one clean wipe up front rather than trying
# to delete a single cell's hyperlink in a loop.
$ws.Range("A2").Hyperlinks.Delete()

$urls = @(
    "https://ae04.alicdn.com/kf/Hce3834858f3c48c783bc3367fbe86e0az.jpg",
    "https://ae04.alicdn.com/kf/H123b859e91564c53aacaaeaebed3d5cd1.jpg",
    "https://ae04.alicdn.com/kf/Hd580e0f5e6284aab8795323ef07dec07J.jpg",
    "https://ae04.alicdn.com/kf/H23b4a341cc554c2c8929017b79efd5fag.jpg",
    "https://ae04.alicdn.com/kf/S6fdd6c761bd94f02b04fbaf386dda424s.jpg",
    "https://ae04.alicdn.com/kf/Hc536e1d0588b4248a79b6fd674abdaf4E.jpg",
    "https://ae04.alicdn.com/kf/Hcefa578dfb3346479f85fa6ae5acb199Z.jpg",
    "https://ae04.alicdn.com/kf/HTB1kXNiAwaTBuNjSszfq6xgfpXaC.jpg"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $urls[$i]
    $ws.Hyperlinks.Add($cell, $urls[$i])
    $cell.Style = "Hyperlink"
}

$ws.Range("A2:A9").Select()
